# "Last fixes & Good night!"
#
# Slide 1, "Subtitle 2" placeholder: the first paragraph's text was built
# from three runs that read "Prepared by" + ":" + ":" (a duplicated colon).
# Fix it by re-splitting the text into "Prepared " + "by" + ":" so the
# paragraph reads "Prepared by:" instead of "Prepared by::".

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Subtitle 2")
$tr  = $shp.TextFrame.TextRange

$full = $tr.Text
$start = $full.IndexOf("Prepared by") + 1   # PowerPoint Characters() is 1-based

# Run 2 (the old, first ":") sits right after "Prepared by" (11 chars).
# Edit it first so the character offsets used for run 1 below are still
# valid (run 1 is edited last, shrinking the text after it).
$run2 = $tr.Characters($start + 11, 1)
$run2.Text = "by"

# Run 1: "Prepared by" -> "Prepared " (keeps the trailing space).
$run1 = $tr.Characters($start, 11)
$run1.Text = "Prepared "
